# Revert commit e236e1f: remove the "weighted value, adjusted for number of
# train cars per locomotive" row (row 37) that was dropped into the
# "BTS NTS Modal Profile Data" sheet, restoring the prior row layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$ws.Rows.Item(37).Delete()

# The "AVLo-passengers" sheet's B5 had been repointed to the new (now
# removed) row; restore the original reference to the un-adjusted weighted
# value in B36 (deleting the row alone would otherwise leave a #REF! here,
# since the formula pointed directly at row 37).
$wsPassengers = $wb.Worksheets.Item("AVLo-passengers")
$wsPassengers.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B36"
